# "Generate Report for Handback"
#
# The handback transform failed for the 0bf7c9c0-...-dad0f36057ac.md file
# in both locales (zh-cn and de-de). Update the per-locale Overview status
# and each locale sheet's Status/Error Detail columns to reflect the
# failure.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhCn = $wb.Worksheets.Item("zh-cn")
$deDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: row 3 (0bf7c9c0-...-dad0f36057ac.md) zh-cn/de-de status cells
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# Per-locale detail sheets: row 3 Status column (C) and new Error Detail (L)
$zhCn.Range("C3").Value = $newStatus
$zhCn.Range("L3").Value = "Handback file name: mofdcp40.uxf is different with handoff file name: 0bf7c9c0-9c19-4c4a-bd91-dad0f36057ac.5991530ca52776e5e1fbbc5a469868305a17b531.zh-cn."

$deDe.Range("C3").Value = $newStatus
$deDe.Range("L3").Value = "Handback file name: mofdcp40.uxf is different with handoff file name: 0bf7c9c0-9c19-4c4a-bd91-dad0f36057ac.5991530ca52776e5e1fbbc5a469868305a17b531.de-de."
